$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Inflammatory-Mac", "Fasl", "Fas", "ECs", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 2, 1, 16.414246, 32.828492, 0.2238644939003909, 0.1639757704305688, 0.9764397232173333, 5.858638339304, 0.03547848010579396, 0.03285174106375043),
  @("Inflammatory-Mac", "Fasl", "Fas", "FAPs", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 3, 1, 8.522109, 25.566327, 0.1162281604801687, 0.1277018197151685, 0.506957538786, 4.562617849074001, 0.01842006477884562, 0.02558443301493019),
  @("Inflammatory-Mac", "Fasl", "Fas", "Inflammatory-Mac", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 3, 1, 8.334395000000001, 25.003185, 0.1136680368163697, 0.1248889691184426, 0.4957909334966668, 4.462118401470001, 0.01801433140464257, 0.02502089219904006),
  @("Inflammatory-Mac", "Fasl", "Fas", "MuSCs", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 2, 1, 3.349162, 6.698324, 0.04567730105424404, 0.03345760866790863, 0.1992327162813333, 1.195396297688, 0.007239027451402953, 0.006703067737899903),
  @("Inflammatory-Mac", "Fasl", "Fas", "Neutrophils", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 3, 1, 34.05638833333334, 102.169165, 0.4644755621620905, 0.5103270520352536, 2.025923724914445, 18.23331352423001, 0.07361098986571545, 0.1022415209714657),
  @("Inflammatory-Mac", "Fasl", "Fas", "Resolving-Mac", 1, 0.3333333333333333, 0.05948733333333334, 0.178462, 0.1584819436421222, 0.2003450935311241, 3, 1, 2.645939, 7.937817, 0.03608644558673645, 0.03964878003265779, 0.1573998552726667, 1.416598697454, 0.005719050035721675, 0.007943438544037791),
  @("MuSCs", "Fasl", "Fas", "ECs", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 2, 1, 16.414246, 32.828492, 0.2238644939003909, 0.1639757704305688, 3.862247462431, 15.448989849724, 0.1403329529732911, 0.08662869848694828),
  @("MuSCs", "Fasl", "Fas", "FAPs", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 3, 1, 8.522109, 25.566327, 0.1162281604801687, 0.1277018197151685, 2.0052394645365, 12.031436787219, 0.07285943695069887, 0.06746510418759792),
  @("MuSCs", "Fasl", "Fas", "Inflammatory-Mac", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 3, 1, 8.334395000000001, 25.003185, 0.1136680368163697, 0.1248889691184426, 1.9610706419075, 11.766423851445, 0.07125458346340323, 0.0659790700888237),
  @("MuSCs", "Fasl", "Fas", "MuSCs", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 2, 1, 3.349162, 6.698324, 0.04567730105424404, 0.03345760866790863, 0.788052794857, 3.152211179428, 0.02863352927974478, 0.01767571566077081),
  @("MuSCs", "Fasl", "Fas", "Neutrophils", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 3, 1, 34.05638833333334, 102.169165, 0.4644755621620905, 0.5103270520352536, 8.013417090250835, 48.08050254150501, 0.2911637575324391, 0.2696067120429495),
  @("MuSCs", "Fasl", "Fas", "Resolving-Mac", 2, 1, 0.2352985, 0.470597, 0.6268656119971068, 0.5283018232478982, 3, 1, 2.645939, 7.937817, 0.03608644558673645, 0.03964878003265779, 0.6225854777915, 3.735512866749001, 0.02262135179752984, 0.02094652278080797),
  @("Neutrophils", "Fasl", "Fas", "ECs", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 2, 1, 16.414246, 32.828492, 0.2238644939003909, 0.1639757704305688, 1.322517685881333, 7.935106115288, 0.04805306082130583, 0.04449533087987007),
  @("Neutrophils", "Fasl", "Fas", "FAPs", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 3, 1, 8.522109, 25.566327, 0.1162281604801687, 0.1277018197151685, 0.6866376849420001, 6.179739164478001, 0.02494865875062417, 0.03465228251264042),
  @("Neutrophils", "Fasl", "Fas", "Inflammatory-Mac", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 3, 1, 8.334395000000001, 25.003185, 0.1136680368163697, 0.1248889691184426, 0.6715133176766668, 6.043619859090001, 0.02439912194832387, 0.03388900683057888),
  @("Neutrophils", "Fasl", "Fas", "MuSCs", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 2, 1, 3.349162, 6.698324, 0.04567730105424404, 0.03345760866790863, 0.2698464478893333, 1.619078687336, 0.009804744323096307, 0.009078825269237917),
  @("Neutrophils", "Fasl", "Fas", "Neutrophils", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 3, 1, 34.05638833333334, 102.169165, 0.4644755621620905, 0.5103270520352536, 2.743968616534445, 24.69571754881001, 0.099700814763936, 0.1384788190208384),
  @("Neutrophils", "Fasl", "Fas", "Resolving-Mac", 1, 0.3333333333333333, 0.08057133333333334, 0.241714, 0.2146524443607711, 0.2713530832209777, 3, 1, 2.645939, 7.937817, 0.03608644558673645, 0.03964878003265779, 0.2131868331486667, 1.918681498338, 0.007746043753484938, 0.01075881870781203)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}